$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "ID" column between EMPLEADO (A) and TAREA (old B) ---
# This shifts the existing TAREA/FECHA/HORA columns from B/C/D to C/D/E while
# preserving their values, shared-string refs and number formats untouched.
$ws.Columns("B").Insert()

# Header for the new column
$ws.Range("B1").Value = "ID"

# --- Fill the ID column for the existing 10 rows (same employees as before) ---
$ws.Range("B2").Value = "@IsaFlores04"
$ws.Range("B3").Value = "@fiamahle"
$ws.Range("B4").Value = "@Ozymandias_96"
$ws.Range("B5").Value = "@DenisBT45"
$ws.Range("B6").Value = "@carlitavot"
$ws.Range("B7").Value = "@diego_roca"
$ws.Range("B8").Value = "@YajaGVargas"
$ws.Range("B9").Value = "@HenryMera"
$ws.Range("B10").Value = "@cesar_sanchez_10"
$ws.Range("B11").Value = "@alexander0266"

# --- Update the two dates/times that changed on the existing rows ---
# ISABEL (row 2): FECHA moves from 45975 to 45980, HORA 0.4375 -> 0.729166...
$ws.Range("D2").Value = 45980
$ws.Range("E2").Value = 0.72916666666666663
# CARLA (row 6): FECHA moves from 45975 to 45994, HORA 0.438194... -> 0.729166...
$ws.Range("D6").Value = 45994
$ws.Range("E6").Value = 0.72916666666666663

# --- Add the new task "LIMPIAR MICROONDAS Y REFRIGERADORA" rows 12-21 ---
$ws.Range("A12").Value = "CARLA"
$ws.Range("B12").Value = "@carlitavot"
$ws.Range("D12").Value = 45975

$ws.Range("A13").Value = "CESAR"
$ws.Range("B13").Value = "@cesar_sanchez_10"
$ws.Range("D13").Value = 46010

$ws.Range("A14").Value = "AHLELY"
$ws.Range("B14").Value = "@fiamahle"
$ws.Range("D14").Value = 45982

$ws.Range("A15").Value = "ISABEL"
$ws.Range("B15").Value = "@IsaFlores04"
$ws.Range("D15").Value = 45987

$ws.Range("A16").Value = "DAVID"
$ws.Range("B16").Value = "@Ozymandias_96"
$ws.Range("D16").Value = 45989

$ws.Range("A17").Value = "DIEGO"
$ws.Range("B17").Value = "@diego_roca"
$ws.Range("D17").Value = 45996

$ws.Range("A18").Value = "ALEXANDER"
$ws.Range("B18").Value = "@alexander0266"
$ws.Range("D18").Value = 46001

$ws.Range("A19").Value = "YAHAIRA"
$ws.Range("B19").Value = "@YajaGVargas"
$ws.Range("D19").Value = 46003

$ws.Range("A20").Value = "DENIS"
$ws.Range("B20").Value = "@DenisBT45"
$ws.Range("D20").Value = 46008

$ws.Range("A21").Value = "HENRY"
$ws.Range("B21").Value = "@HenryMera"
$ws.Range("D21").Value = 46010

$ws.Range("C12:C21").Value = "LIMPIAR MICROONDAS Y REFRIGERADORA"
$ws.Range("E12:E21").Value = 0.66666666666666663

# Give the new FECHA/HORA cells (rows 12-21) the same number formats used by
# the existing FECHA (D) / HORA (E) columns, reusing the existing styles
# instead of creating new ones.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D12:D21").PasteSpecial(-4122) | Out-Null

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E12:E21").PasteSpecial(-4122) | Out-Null

# --- Selection shown in the saved file ---
$ws.Range("E13:E21").Select()
